# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches from the deck's custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2. The presentation theme (theme1.xml, used by the slide master) and the
#     notes/handout theme (theme2.xml) had their contents swapped: the main
#     deck theme becomes the default "Office Theme" color palette, while the
#     notes theme becomes the "Integral" palette that used to be the main
#     theme.  This runtime only exposes the slide-master theme through the
#     object model (NotesMaster/HandoutMaster alias back to the same
#     SlideMaster/theme1.xml here), so the reachable, faithful part of that
#     swap is re-pointing the live deck's theme colors at the "Office Theme"
#     palette via ThemeColorScheme.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $candidate = $slide6.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
$tableShape.Table.ApplyStyle("{C832FF3D-E6A4-436B-B55A-AA55D2F68756}")

# --- 2. Theme color swap (Integral -> Office Theme) -----------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0x000000  # dk1
$tcs.Colors(2).RGB  = 0xFFFFFF  # lt1
$tcs.Colors(3).RGB  = 0x6A5444  # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED  # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244  # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70  # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305  # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95  # folHlink -> 954F72
